$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.811.28'
$ws.Range("E2").Value = '  +1.12%  '
$ws.Range("D3").Value = '1.885.76'
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("D4").Value = '''1.018'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +1.68%  '
$ws.Range("D5").Value = '''334.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.68%  '
$ws.Range("E6").Value = '  +1.52%  '
$ws.Range("D7").Value = '''0.4675'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("D8").Value = '''0.3912'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.66%  '
$ws.Range("D9").Value = '''47.42'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.31%  '
$ws.Range("D10").Value = '''0.08041'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.18%  '
$ws.Range("D11").Value = '''1.014'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("D12").Value = '''21.91'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").Value = '1.913.59'
$ws.Range("E13").Value = '  +2.50%  '
$ws.Range("D14").Value = '''5.954'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.06%  '
$ws.Range("D15").Value = '''7.078'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = '''1.020'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.63%  '
$ws.Range("D17").Value = '''0.06749'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.87%  '
$ws.Range("D18").Value = '''87.24'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("D19").Value = '''0.00001048'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.10%  '
$ws.Range("D20").Value = '''17.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.82%  '
$ws.Range("D21").Value = '''1.016'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.47%  '
$ws.Range("D22").Value = '27.858.25'
$ws.Range("E22").Value = '  +1.25%  '
$ws.Range("D23").Value = '''5.493'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.38%  '
$ws.Range("D24").Value = '''10.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.68%  '
$ws.Range("E25").Value = '  +1.63%  '
$ws.Range("D26").Value = '2.138.45'
$ws.Range("E26").Value = '  +2.41%  '
$ws.Range("D27").Value = '''159.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.67%  '
$ws.Range("D28").Value = '''20.05'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.11%  '
$ws.Range("D29").Value = '''2.078'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.36%  '
$ws.Range("D30").Value = '''5.459'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.98%  '
$ws.Range("D31").Value = '''121.80'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.50%  '
$ws.Range("D32").Value = '''0.9670'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = '''0.09477'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("D34").Value = '''3.645'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.17%  '
$ws.Range("D35").Value = '''1.408'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.31%  '
$ws.Range("D36").Value = '''5.345'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.79%  '
$ws.Range("D37").Value = '''0.06115'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("D38").Value = '''0.02254'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("D39").Value = '''1.211'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.75%  '
$ws.Range("D40").Value = '''0.5976'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.33%  '
$ws.Range("D41").Value = '''7.993'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.94%  '
$ws.Range("D42").Value = '''0.1884'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.52%  '
$ws.Range("D43").Value = '''10.26'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.49%  '
$ws.Range("D44").Value = '''1.265'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D45").Value = '''0.5671'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("D46").Value = '''12.18'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.40%  '
$ws.Range("D47").Value = '''3.403'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.24%  '
$ws.Range("D48").Value = '''1.922'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.01%  '
$ws.Range("D49").Value = '''0.06918'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.87%  '
$ws.Range("D50").Value = '''113.67'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.13%  '
$ws.Range("D51").Value = '''1.069'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.46%  '
